$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the c27a0754-... row on both the zh-cn and de-de report sheets,
# as part of regenerating the handback report.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-13 08:34:14"
$wsZhCn.Range("G3").Value = "2016-01-13 08:35:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-13 08:34:38"
$wsDeDe.Range("G3").Value = "2016-01-13 08:36:04"
